# Update the model-comparison data (umwm vs planetwaves) with refreshed RMS values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 0.030421299999999998
$ws.Range("B5").Value = 0.081498347222804995
$ws.Range("C5").Value = 0.14451800000000001
$ws.Range("B6").Value = 0.18736366927623699
$ws.Range("C6").Value = 0.26693800000000001
$ws.Range("B7").Value = 0.319779723882675
$ws.Range("C7").Value = 0.43732599999999999
$ws.Range("B8").Value = 0.48304304480552601
$ws.Range("C8").Value = 0.64294499999999999
$ws.Range("B9").Value = 0.67869186401367099
$ws.Range("C9").Value = 0.87729599999999996
$ws.Range("B10").Value = 0.90832310914993197
$ws.Range("C10").Value = 1.14096
$ws.Range("B11").Value = 1.1727869510650599
$ws.Range("C11").Value = 1.4317200000000001
$ws.Range("B12").Value = 1.47175681591033
$ws.Range("C12").Value = 1.7474400000000001
$ws.Range("B13").Value = 1.80442786216735
$ws.Range("C13").Value = 2.0729600000000001
$ws.Range("B14").Value = 2.1610410213470401
$ws.Range("C14").Value = 2.42971
$ws.Range("B15").Value = 2.5583541393279998
$ws.Range("C15").Value = 2.7871100000000002
$ws.Range("B16").Value = 2.98986363410949
$ws.Range("C16").Value = 3.1803699999999999
$ws.Range("B17").Value = 3.4347374439239502
$ws.Range("C17").Value = 3.5927699999999998
$ws.Range("B18").Value = 3.93021535873413
$ws.Range("C18").Value = 3.9786000000000001
$ws.Range("B19").Value = 4.4217596054077104
$ws.Range("C19").Value = 4.3825599999999998
$ws.Range("B20").Value = 4.9297094345092702
$ws.Range("C20").Value = 4.8433299999999999
$ws.Range("B21").Value = 5.4979319572448704
$ws.Range("C21").Value = 5.2761300000000002
$ws.Range("B22").Value = 6.0374197959899902
$ws.Range("C22").Value = 5.77555

# Select column D (to mirror the saved cursor/selection state from the edit session)
$ws.Range("D1:D1048576").Select()
